$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data table (rows 2-58) holds one price record per row. A new weekly
# record (two quality grades: "Primera" and "Segunda") for date 2022-07-27
# (serial 44769) is inserted right after the existing 2022-07-15 (44757)
# entries, pushing every following record down by two rows.
$ws.Rows("52:53").Insert()

# New row 52: "Primera" quality record for 44769
$ws.Range("A52").Value = 11
$ws.Range("B52").Value = "Vega Monumental Concepción"
$ws.Range("C52").Value = "Bíobío"
$ws.Range("D52").Value = 44769
$ws.Range("D52").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E52").Value = 8
$ws.Range("F52").Value = 100112037
$ws.Range("G52").Value = "Cebollín"
$ws.Range("H52").Value = "Sin especificar"
$ws.Range("I52").Value = "Primera"
$ws.Range("J52").Value = 100
$ws.Range("K52").Value = 700
$ws.Range("L52").Value = 800
$ws.Range("M52").Value = 750
$ws.Range("N52").Value = "$/paquete 6 unidades"
$ws.Range("O52").Value = "Región de Ñuble"
$ws.Range("P52").Value = 125
$ws.Range("Q52").Value = 6
$ws.Range("R52").Value = "Hortaliza"

# New row 53: "Segunda" quality record for 44769
$ws.Range("A53").Value = 11
$ws.Range("B53").Value = "Vega Monumental Concepción"
$ws.Range("C53").Value = "Bíobío"
$ws.Range("D53").Value = 44769
$ws.Range("D53").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E53").Value = 8
$ws.Range("F53").Value = 100112037
$ws.Range("G53").Value = "Cebollín"
$ws.Range("H53").Value = "Sin especificar"
$ws.Range("I53").Value = "Segunda"
$ws.Range("J53").Value = 50
$ws.Range("K53").Value = 600
$ws.Range("L53").Value = 600
$ws.Range("M53").Value = 600
$ws.Range("N53").Value = "$/paquete 6 unidades"
$ws.Range("O53").Value = "Región de Ñuble"
$ws.Range("P53").Value = 100
$ws.Range("Q53").Value = 6
$ws.Range("R53").Value = "Hortaliza"
